# Add a new "Audio" topic section to the glossary table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New glossary rows (40-44) ---
# Fill Term (column B) first for all new rows.
$ws.Cells.Item(40, 2).Value = "Amplifier"
$ws.Cells.Item(41, 2).Value = "Bus"
$ws.Cells.Item(42, 2).Value = "Channel aka Strip"
$ws.Cells.Item(43, 2).Value = "Direct Box"
$ws.Cells.Item(44, 2).Value = "Scene"

# Then fill Topic (column A) for all the new rows.
$ws.Cells.Item(40, 1).Value = "Audio"
$ws.Cells.Item(41, 1).Value = "Audio"
$ws.Cells.Item(42, 1).Value = "Audio"
$ws.Cells.Item(43, 1).Value = "Audio"
$ws.Cells.Item(44, 1).Value = "Audio"

# Finally fill Definition (column D) for all the new rows.
$ws.Cells.Item(40, 4).Value = "An audio component that take low-level inputs, from a console, and amplifies them to the power level required to drive a speaker. "
$ws.Cells.Item(41, 4).Value = "Refers to a logical structure within the console where multiple channels (sources) can contribute an audio signal."
$ws.Cells.Item(42, 4).Value = "Refers to an input in the console. Typically this is a microphone input or an instrument."
$ws.Cells.Item(43, 4).Value = "A device, usually a box, which can interface an instrument, like a electronic keyboard, to an audio console."
$ws.Cells.Item(44, 4).Value = "A collection of setting which can be saved and recalled as a group."

# --- Grow the glossary table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D53"))

# --- Apply wrap-text formatting to the (currently empty) trailing column-D cells ---
$ws.Range("D45:D53").WrapText = $true

# --- Update the visible selection/scroll position to match where editing left off ---
$ws.Range("D44").Select()
